$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 4) matching the existing table's layout/format.
$ws.Range("A4").Value = "'20240223"
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B4").Value = 174
$ws.Range("C4").Value = 674
$ws.Range("D4").Value = 345
$ws.Range("E4").Value = 86
$ws.Range("F4").Value = 45
$ws.Range("G4").Value = 152
